$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Helper cell (far off the used range) used purely to push a TEXT-typed
# value into date-look-alike cells without Excel's auto date-detection
# kicking in (which would otherwise convert "06-01-2025" to a date serial
# and swap in a brand-new number-format style).
$helper = $ws1.Range("ZZ1")
$helper.NumberFormat = "@"

function Set-TextValue {
    param($range, [string]$text)
    $helper.Value = $text
    $helper.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# --- Sheet1 (wfCreateEntity_OOFS_MultiSuccessCase row) ---
Set-TextValue $ws1.Range("O2") "06-01-2025"
$ws1.Range("Q2").Value = "09-01-2025 05:00:00 PM"
Set-TextValue $ws1.Range("AD2") "06-01-2025"
$ws1.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 7:50 PM"

# --- Sheet2 (wfCreateEntity_OOFS_SuccessUpdate row) ---
Set-TextValue $ws2.Range("O2") "06-01-2025"
$ws2.Range("Q2").Value = "09-01-2025 05:00:00 PM"
Set-TextValue $ws2.Range("AD2") "06-01-2025"
$ws2.Range("BB2").Value = "CT: Mon, Jan 06, 2025 at 8:12 PM"

$helper.Clear() | Out-Null
$excel.CutCopyMode = $false
